$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for rows 2-7, columns G (7) through T (20)
$values = @{
    2 = @{ G=13.10570166666667;  H=39.317105;               I=0.004355939447658156; J=0.004355939447658156;
           K=3; L=1; M=0.110028; N=0.330084;
           Q=1.44199414298;      R=12.97794728682;           S=0.004355939447658156; T=0.004355939447658156 }
    3 = @{ I=0.001077974419837672; J=0.001077974419837672;
           K=3; L=1; M=0.110028; N=0.330084;
           Q=0.35685362902;       R=3.21168266118;            S=0.001077974419837672; T=0.001077974419837672 }
    4 = @{ G=608.3979493333333;  H=1825.193848;             I=0.2022131050118309;  J=0.202213105011831;
           K=3; L=1; M=0.110028; N=0.330084;
           Q=66.94080956924799;  R=602.4672861232319;        S=0.2022131050118309;  T=0.202213105011831 }
    5 = @{ G=3.045399333333334; H=9.136198;                 I=0.001012198768699159; J=0.001012198768699159;
           K=3; L=1; M=0.110028; N=0.330084;
           Q=0.3350791978480001; R=3.015712780632;           S=0.001012198768699159; T=0.001012198768699159 }
    6 = @{ G=255.8380176666667; H=767.514053;               I=0.08503283087843555; J=0.08503283087843555;
           K=3; L=1; M=0.110028; N=0.330084;
           Q=28.149345407828;    R=253.344108670452;         S=0.08503283087843555; T=0.08503283087843555 }
    7 = @{ G=2125.066569;       H=6375.199707000001;        I=0.7063079514735385;  J=0.7063079514735385;
           K=3; L=1; M=0.110028; N=0.330084;
           Q=233.816824453932;   R=2104.351420085388;        S=0.7063079514735385;  T=0.7063079514735385 }
}

foreach ($row in $values.Keys) {
    $cols = $values[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
